# Improved Allies: add two new spell rows (friend/attack, friend/buff) to the
# spellDictionary sheet, inserted above the existing "illyia"/style-list block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 17, pushing the existing rows 17-29 down to 19-31.
$ws.Rows("17:18").Insert()

# --- Row 18: "illyia" / the ultimate sword (friend/attack) ---
# Filled in the same order the original author appears to have used: A, C, D-K, then B.
$ws.Cells.Item(18, 1).Value = "illyia"
$ws.Cells.Item(18, 3).Value = "the ultimate sword"
$ws.Cells.Item(18, 4).Value = 999
$ws.Cells.Item(18, 5).Value = 30
$ws.Cells.Item(18, 6).Value = 100
$ws.Cells.Item(18, 7).Value = 75
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = "MT"
$ws.Cells.Item(18, 10).Value = "ADSHEUFIB"
$ws.Cells.Item(18, 11).Value = 100
$ws.Cells.Item(18, 2).Value = "friend/attack"

# --- Row 17: "dahlia" / Super Buff! (friend/buff) ---
$ws.Cells.Item(17, 1).Value = "dahlia"
$ws.Cells.Item(17, 2).Value = "friend/buff"
$ws.Cells.Item(17, 3).Value = "Super Buff!"
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 25
$ws.Cells.Item(17, 6).Value = 100
$ws.Cells.Item(17, 7).Value = 10
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = "m"
$ws.Cells.Item(17, 10).Value = "adsheufib"
$ws.Cells.Item(17, 11).Value = 100

# Update the selected cell to match the author's final cursor position.
$null = $ws.Range("K18").Select()
